# ValueSet-wh-payer-education-level: metadata refresh (FHIR IG Alvearie build).
#
#   Version     5.0.0  -> 6.0.0
#   Date        2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
#   Publisher   (blank) -> "Alvearie Team"
#
# The old "Contact" / "No display for ContactDetail" row is replaced by a
# "Jurisdiction" / "United States of America" row, and the duplicate
# "Contact" row that used to sit right below it is removed outright, so the
# Metadata sheet shrinks from 15 rows to 14.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B3").Value = "6.0.0"
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"
$ws.Range("B9").Value = "Alvearie Team"

$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Delete the now-redundant second "Contact" row; rows 12-15 shift up to
# become rows 11-14.
$ws.Rows.Item(11).Delete()
